$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$epoch = [DateTime]::FromOADate(25569)
$lastRow = $ws.Cells.Item(1,1).End(-4121).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $ms = $cell.Value()
    if ($ms -ne $null) {
        $sec = $ms / 1000
        $dt = $epoch.AddSeconds($sec)
        $cell.Value = $dt.ToString("yyyy-MM-dd HH:mm:ss")
    }
}
